$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "31.220.49"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "1.942.72"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "242.34"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "0.2914"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "0.06795"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "20.28"
$ws.Range("E10").Value = "  +7.98%  "
$ws.Range("D11").Value = "104.44"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07850"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "1.957.30"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "5.308"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "0.6969"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "297.35"
$ws.Range("E16").Value = "  +11.83%  "
$ws.Range("D17").Value = "31.223.36"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "2.208.95"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "0.000007618"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "5.575"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "6.446"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "9.569"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "168.94"
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").Value = "19.86"
$ws.Range("E27").Value = "  +4.04%  "
$ws.Range("D28").Value = "2.106"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "1.403"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "4.638"
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").Value = "1.537"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "4.355"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "0.04833"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").Value = "0.7413"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "1.138"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "6.909"
$ws.Range("E39").Value = "  +8.75%  "
$ws.Range("D40").Value = "2.641"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "76.86"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "2.034"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "0.8739"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "0.4383"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "106.18"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "1.025.35"
$ws.Range("E47").Value = "  +7.54%  "
$ws.Range("D48").Value = "7.593"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "9.232"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "35.23"
$ws.Range("E51").Value = "  -0.01%  "
